# Second commit: remove the "bobby2" row (row 3, with its hyperlink) and
# add a new "Country" column (F) populated with "India" for the remaining
# data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink anchored on B3 before the row shifts up, so it
# doesn't linger as an orphaned hyperlink entry.
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -eq '$B$3') {
        $h.Delete()
    }
}

# Delete the entire third row (bobby2 / test2@mail.com / testphone2 / Female / Monday).
$ws.Rows.Item(3).Delete()

# Add the new "Country" column header and value.
$ws.Range("F1").Value = "Country"
$ws.Range("F2").Value = "India"
